$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the standalone "Meta description" paragraph that currently sits
#    right under the H1 title. It is a 3-run paragraph:
#      <w:r/>  <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#      <w:r><w:t>: Get into the Everglades ... jackpots!</w:t></w:r>
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Meta description")) {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) The final paragraph in the document (the image "Prompt: ..." paragraph)
#    gets a brand-new bold paragraph inserted right before it, carrying the
#    text that used to live in the "Meta description" run, and its own text
#    is swapped for the (now un-prefixed) meta-description sentence.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)

# Create the paragraph break first so the new paragraph has no inherited
# pPr/rPr baggage, then stamp its exact run structure via InsertXML.
$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:r/>' +
              '<w:r><w:rPr><w:b/></w:rPr>' +
              '<w:t>Play 15 Armadillos Slot Machine | Free Spins and Jackpots</w:t>' +
              '</w:r></w:p>'
$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3) Swap the text of the (still) last paragraph from the image prompt to
#    the meta-description sentence, keeping its italic run formatting as-is.
# ---------------------------------------------------------------------------
$oldPromptText = "Prompt: Create a cartoon-style feature image for the game " +
    [char]34 + "15 Armadillos" + [char]34 + " that features a happy Maya warrior with glasses. " +
    "For the feature image of " + [char]34 + "15 Armadillos" + [char]34 + ", let's have a cartoon-style design featuring a happy Maya warrior with glasses. " +
    "The warrior can be depicted wearing a headdress made of colorful feathers, with intricate designs on their face and body. " +
    "They can be holding a staff or weapon made of stone or wood, with a happy expression on their face. " +
    "In the background, we can see the Everglades National Park with its lush greenery and animals like alligators and otters. " +
    "The image can be bright and colorful to reflect the fun and adventurous nature of the game."

$newPromptText = "Get into the Everglades for a chance to win big with the 15 Armadillos slot machine! " +
    "Play now and enjoy free spins, pick accumulation, and five jackpots!"

$d.Content.Find.Execute($oldPromptText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newPromptText, 2)
